# Update Name of Algo
# Apply updated numeric values (re-run of RandomForest imputation) to result_data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = -11.6878
$ws.Range("B3").Value  = 6.096099999999987
$ws.Range("C5").Value  = -14.22670000000001
$ws.Range("D5").Value  = -8.810699999999994
$ws.Range("E7").Value  = 12.0711
$ws.Range("D9").Value  = -8.673700000000007
$ws.Range("D11").Value = -8.356300000000006
$ws.Range("E11").Value = 13.4511
$ws.Range("B14").Value = 9.202200000000001
$ws.Range("B16").Value = 9.193100000000001
$ws.Range("C16").Value = -11.22349999999999
$ws.Range("D17").Value = -8.735500000000002
$ws.Range("E19").Value = 13.36459999999999
$ws.Range("B21").Value = 5.528399999999993
$ws.Range("D21").Value = -7.867500000000007
$ws.Range("E21").Value = 13.7941
$ws.Range("B23").Value = 5.950299999999998
$ws.Range("B25").Value = 5.655499999999995
